# The deck originally had 3 slides:
#   1. "高性能、高可用场景带来的挑战"  (sldId 257)
#   2. "后端架构展示图" backend architecture diagram (sldId 256)
#   3. "我们还有很多工作…"            (sldId 259)
#
# The edit keeps only the backend architecture diagram slide and drops
# the other two (along with the notes page that belonged to slide 1).
# Removing the higher-index slide first keeps slide 1's index stable.
$p = $ppt.ActivePresentation

$p.Slides.Item(3).Delete()
$p.Slides.Item(1).Delete()
